$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.590.82"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "3.595.97"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.489"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "4.206.95"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.638.88"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "66.678.89"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.622"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "3.737.76"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "3.593.15"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "176.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0857"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.897"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.44%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.28%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.02%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.237"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
